$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Study")
$ws.Name = "isa_study"
$ws.Activate()
